$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the FFT Median RF/VM/VL columns (K:M), which also drops the
# now-unused "FFT Median RF", "FFT Median VM", "FFT Median VL" and "NA"
# shared strings and shrinks the sheet dimension to A1:J56.
$ws.Columns("K:M").Delete() | Out-Null

# Update the 5th measurement row for each subject (t1..t5) with corrected
# values for columns C through J.
function Set-RowValues($row, $values) {
    $cols = @("C", "D", "E", "F", "G", "H", "I", "J")
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value = $values[$i]
    }
}

Set-RowValues 6  @(132.75931, 571.9058, 2.8665419999999999, 3.1026600000000002, 3.1727729999999998, 82.849199999999996, 78.839399999999998, 69.237499999999997)
Set-RowValues 17 @(106.16804, 715.32790999999997, 2.8588399999999998, 3.1006, 3.1993200000000002, 76.352689999999996, 74.52534, 65.849199999999996)
Set-RowValues 28 @(117.773, 661.58789999999999, 2.61795, 2.8025129999999998, 2.8969040000000001, 85.707018000000005, 78.628860000000003, 66.129790200000002)
Set-RowValues 39 @(116.60128, 647.10428100000001, 2.5696189999999999, 2.8909159999999998, 3.0233143999999998, 86.187200000000004, 76.937200000000004, 66.650189999999995)
Set-RowValues 50 @(113.25248999999999, 687.14979900000003, 2.6093860000000002, 2.986866, 2.9538850000000001, 85.312619999999995, 77.272715000000005, 66.448560999999998)

# Restore the original selection/view state.
$ws.Range("D6").Select() | Out-Null
